# This script applies updated Markov transition-matrix probabilities to
# Sheet1 (team_specific_matrix/UMBC_A). The underlying simulation was
# re-run with more games and faster simulate-game logic, which changed
# the empirical transition frequencies recorded in B2:S19. Each row
# B:S still sums to 1 (it is a row-stochastic transition matrix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2071713147410359
$ws.Range("C2").Value = 0.5099601593625498
$ws.Range("J2").Value = 0.00796812749003984
$ws.Range("P2").Value = 0.1633466135458167
$ws.Range("S2").Value = 0.1115537848605578
$ws.Range("B3").Value = 0.0234375
$ws.Range("C3").Value = 0.0390625
$ws.Range("J3").Value = 0.0234375
$ws.Range("P3").Value = 0.7265625
$ws.Range("S3").Value = 0.1875
$ws.Range("J4").Value = 0.0975609756097561
$ws.Range("P4").Value = 0.5365853658536586
$ws.Range("S4").Value = 0.3658536585365854
$ws.Range("B6").Value = 0.06306306306306306
$ws.Range("D6").Value = 0.01351351351351351
$ws.Range("F6").Value = 0.04504504504504504
$ws.Range("J6").Value = 0.2567567567567567
$ws.Range("O6").Value = 0.009009009009009009
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("R6").Value = 0.04054054054054054
$ws.Range("S6").Value = 0.4054054054054054
$ws.Range("B7").Value = 0.09743589743589744
$ws.Range("D7").Value = 0.01538461538461539
$ws.Range("F7").Value = 0.04102564102564103
$ws.Range("J7").Value = 0.1333333333333333
$ws.Range("O7").Value = 0.01025641025641026
$ws.Range("Q7").Value = 0.1897435897435897
$ws.Range("R7").Value = 0.09230769230769231
$ws.Range("S7").Value = 0.4205128205128205
$ws.Range("B8").Value = 0.06702898550724638
$ws.Range("D8").Value = 0.02355072463768116
$ws.Range("F8").Value = 0.05253623188405797
$ws.Range("J8").Value = 0.1141304347826087
$ws.Range("O8").Value = 0.01630434782608696
$ws.Range("Q8").Value = 0.2318840579710145
$ws.Range("R8").Value = 0.09057971014492754
$ws.Range("S8").Value = 0.4039855072463768
$ws.Range("B9").Value = 0.1161616161616162
$ws.Range("D9").Value = 0.04040404040404041
$ws.Range("F9").Value = 0.06565656565656566
$ws.Range("J9").Value = 0.0505050505050505
$ws.Range("O9").Value = 0.01515151515151515
$ws.Range("Q9").Value = 0.202020202020202
$ws.Range("R9").Value = 0.0707070707070707
$ws.Range("S9").Value = 0.4393939393939394
$ws.Range("B10").Value = 0.0796812749003984
$ws.Range("D10").Value = 0.01274900398406375
$ws.Range("E10").Value = 0.002390438247011952
$ws.Range("F10").Value = 0.07569721115537849
$ws.Range("J10").Value = 0.09880478087649402
$ws.Range("O10").Value = 0.01035856573705179
$ws.Range("Q10").Value = 0.2310756972111554
$ws.Range("R10").Value = 0.100398406374502
$ws.Range("S10").Value = 0.3888446215139442
$ws.Range("G11").Value = 0.1684981684981685
$ws.Range("J11").Value = 0.06227106227106227
$ws.Range("K11").Value = 0.1941391941391941
$ws.Range("L11").Value = 0.5641025641025641
$ws.Range("S11").Value = 0.01098901098901099
$ws.Range("G12").Value = 0.7423312883435583
$ws.Range("J12").Value = 0.1840490797546012
$ws.Range("K12").Value = 0.006134969325153374
$ws.Range("L12").Value = 0.03680981595092025
$ws.Range("S12").Value = 0.03067484662576687
$ws.Range("F13").Value = 0.02127659574468085
$ws.Range("G13").Value = 0.6595744680851063
$ws.Range("J13").Value = 0.2978723404255319
$ws.Range("S13").Value = 0.02127659574468085
$ws.Range("G14").Value = 0.75
$ws.Range("J14").Value = 0.25
$ws.Range("F15").Value = 0.01176470588235294
$ws.Range("H15").Value = 0.1764705882352941
$ws.Range("I15").Value = 0.07647058823529412
$ws.Range("J15").Value = 0.3235294117647059
$ws.Range("K15").Value = 0.07647058823529412
$ws.Range("M15").Value = 0.01176470588235294
$ws.Range("O15").Value = 0.05294117647058823
$ws.Range("S15").Value = 0.2705882352941176
$ws.Range("F16").Value = 0.02597402597402598
$ws.Range("H16").Value = 0.2012987012987013
$ws.Range("I16").Value = 0.08441558441558442
$ws.Range("J16").Value = 0.3571428571428572
$ws.Range("K16").Value = 0.1168831168831169
$ws.Range("M16").Value = 0.01298701298701299
$ws.Range("O16").Value = 0.05194805194805195
$ws.Range("S16").Value = 0.1493506493506493
$ws.Range("F17").Value = 0.02087286527514232
$ws.Range("H17").Value = 0.2030360531309298
$ws.Range("I17").Value = 0.09677419354838709
$ws.Range("J17").Value = 0.396584440227704
$ws.Range("K17").Value = 0.09297912713472485
$ws.Range("M17").Value = 0.01518026565464896
$ws.Range("N17").Value = 0.003795066413662239
$ws.Range("O17").Value = 0.04364326375711575
$ws.Range("S17").Value = 0.127134724857685
$ws.Range("F18").Value = 0.01851851851851852
$ws.Range("H18").Value = 0.2453703703703704
$ws.Range("I18").Value = 0.06944444444444445
$ws.Range("J18").Value = 0.3703703703703703
$ws.Range("K18").Value = 0.1111111111111111
$ws.Range("M18").Value = 0.02314814814814815
$ws.Range("O18").Value = 0.05555555555555555
$ws.Range("S18").Value = 0.1064814814814815
$ws.Range("F19").Value = 0.01316752011704462
$ws.Range("H19").Value = 0.246525237746891
$ws.Range("I19").Value = 0.07827359180687637
$ws.Range("J19").Value = 0.3803950256035113
$ws.Range("K19").Value = 0.08339429407461595
$ws.Range("M19").Value = 0.02194586686174104
$ws.Range("N19").Value = 0.002194586686174104
$ws.Range("O19").Value = 0.04901243599122165
$ws.Range("S19").Value = 0.1250914411119239
